$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "methods" (shipping method) column for every data row:
# "Ground" -> "Ground Shipping"
$ws.Range("V2:V27").Value = "Ground Shipping"

# Update the active selection to match the author's final cursor position
$ws.Range("V33").Select()
